$d = $word.ActiveDocument

# Locate the paragraph that holds the " m:'doc.html'.fromHTMLURI() " field
# (fldChar begin / instrText* / fldChar end) and rewrite it in place as
# literal text runs "{" m "," ":" "'" "doc.html" "'.fromHTMLURI()" "}" while
# keeping the _GoBack bookmark pair exactly where it was, between
# "doc.html" and "'.fromHTMLURI()".

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
    }
}

$rsidR = $target.Range.ParagraphFormat.Parent.Range.Information(1)

$p = $target
$r = $p.Range

$apos = [char]39

$inner = ""
$inner += "<w:r><w:t>{</w:t></w:r>"
$inner += "<w:r><w:t>m</w:t></w:r>"
$inner += "<w:r><w:t>:</w:t></w:r>"
$inner += "<w:r><w:t>$apos</w:t></w:r>"
$inner += "<w:r><w:t>doc.html</w:t></w:r>"
$inner += "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>"
$inner += "<w:bookmarkEnd w:id=`"0`"/>"
$inner += "<w:r><w:t>$apos.fromHTMLURI()</w:t></w:r>"
$inner += "<w:r><w:t xml:space=`"preserve`">}</w:t></w:r>"

$pOpen = "<w:p w:rsidR=`"00C52979`" w:rsidRDefault=`"00C52979`" w:rsidP=`"00F5495F`">"

$xml = '<?xml version="1.0" standalone="yes"?>' +
       '<?mso-application progid="Word.Document"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' + $pOpen + $inner + '</w:p></w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
